$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 previously held the text "RESIDENT-1"; change it to the numeric value 1001.
$ws.Range("B2").Value = 1001

# Match the selection highlighted in the saved workbook (whole used range).
$ws.Range("A1:I2").Select() | Out-Null
